$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.250.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "'2.789.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'346.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("D6").Value = "'115.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "'42.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").Value = "'0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "'20.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D14").Value = "'7.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "'3.228.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "'2.770.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "'0.889"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'52.109.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "'3.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.76%  "
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").Value = "'13.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").Value = "'0.0₃0977"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'269.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'69.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'2.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.75%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "'0.139"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'34.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").Value = "'50.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'0.0448"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +28.07%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'0.0824"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'4.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("E39").Value = "  -4.71%  "
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").Value = "'2.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.88%  "
$ws.Range("D42").Value = "'127.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "'23.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "'2.062.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "'0.963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.50%  "
$ws.Range("D50").Value = "'5.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -1.74%  "
